# Auto-generated edit script: updates currentAveragePrice / Leve price & profit
# columns (H..N) for specific rows across multiple sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = $null
$ws.Range("H99").Value = 257.5
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 257.5
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 772.5
$ws.Range("N99").Value = -3768.5
$ws.Range("H111").Value = 3127098.2
$ws.Range("I111").Value = 3354.125
$ws.Range("J111").Value = 4168346.2
$ws.Range("K111").Value = 10062.375
$ws.Range("L111").Value = 12505038.6
$ws.Range("M111").Value = -6995.375
$ws.Range("N111").Value = -12511172.6
$ws.Range("H113").Value = 2708.7727
$ws.Range("I113").Value = 2720
$ws.Range("J113").Value = 2699.4167
$ws.Range("K113").Value = 2720
$ws.Range("L113").Value = 2699.4167
$ws.Range("M113").Value = 534
$ws.Range("N113").Value = -9207.4167
$ws.Range("H125").Value = 2497
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2497
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 22473
$ws.Range("N125").Value = -27393
$ws.Range("M125").Value = $null
$ws.Range("H138").Value = 3630.4
$ws.Range("I138").Value = 3009.4
$ws.Range("J138").Value = 3754.6
$ws.Range("K138").Value = 9028.200000000001
$ws.Range("L138").Value = 11263.8
$ws.Range("M138").Value = -3888.200000000001
$ws.Range("N138").Value = -21543.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3930.1428
$ws.Range("I45").Value = 3502.2
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 3502.2
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -3125.2
$ws.Range("N45").Value = -5754

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3006.5
$ws.Range("I20").Value = 2738.5881
$ws.Range("J20").Value = 3512.5557
$ws.Range("K20").Value = 2738.5881
$ws.Range("L20").Value = 3512.5557
$ws.Range("M20").Value = -2491.5881
$ws.Range("N20").Value = -4006.5557
$ws.Range("H86").Value = 369615.1
$ws.Range("I86").Value = 6396.4
$ws.Range("J86").Value = 672297.3
$ws.Range("K86").Value = 6396.4
$ws.Range("L86").Value = 672297.3
$ws.Range("M86").Value = -5273.4
$ws.Range("N86").Value = -674543.3
$ws.Range("H89").Value = 369615.1
$ws.Range("I89").Value = 6396.4
$ws.Range("J89").Value = 672297.3
$ws.Range("K89").Value = 31982
$ws.Range("L89").Value = 3361486.5
$ws.Range("M89").Value = -26366
$ws.Range("N89").Value = -3372718.5
$ws.Range("H94").Value = 3164.9167
$ws.Range("I94").Value = 2569.8572
$ws.Range("J94").Value = 3998
$ws.Range("K94").Value = 2569.8572
$ws.Range("L94").Value = 3998
$ws.Range("M94").Value = -2118.8572
$ws.Range("N94").Value = -4900
$ws.Range("H99").Value = 2355.4666
$ws.Range("I99").Value = 2459.7273
$ws.Range("J99").Value = 2068.75
$ws.Range("K99").Value = 2459.7273
$ws.Range("L99").Value = 2068.75
$ws.Range("M99").Value = -961.7273
$ws.Range("N99").Value = -5064.75
$ws.Range("H105").Value = 200051800
$ws.Range("I105").Value = 333417200
$ws.Range("J105").Value = 3725
$ws.Range("K105").Value = 333417200
$ws.Range("L105").Value = 3725
$ws.Range("M105").Value = -333415453
$ws.Range("N105").Value = -7219
$ws.Range("H107").Value = 3497.7896
$ws.Range("I107").Value = 2900.4119
$ws.Range("J107").Value = 8575.5
$ws.Range("K107").Value = 2900.4119
$ws.Range("L107").Value = 8575.5
$ws.Range("M107").Value = -980.4119000000001
$ws.Range("N107").Value = -12415.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 300
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 300
$ws.Range("N6").Value = -526
$ws.Range("H7").Value = 233.05263
$ws.Range("I7").Value = 67.333336
$ws.Range("J7").Value = 382.2
$ws.Range("K7").Value = 67.333336
$ws.Range("L7").Value = 382.2
$ws.Range("M7").Value = 45.666664
$ws.Range("N7").Value = -608.2
$ws.Range("H11").Value = 2750
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 4500
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 4500
$ws.Range("M11").Value = -860
$ws.Range("N11").Value = -4780
$ws.Range("H12").Value = 275000
$ws.Range("I12").Value = 275000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 275000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -274830
$ws.Range("H17").Value = 10008
$ws.Range("I17").Value = 10008
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 10008
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -9834
$ws.Range("N17").Value = $null
$ws.Range("H19").Value = 898.2857
$ws.Range("I19").Value = 464.66666
$ws.Range("J19").Value = 3500
$ws.Range("K19").Value = 464.66666
$ws.Range("L19").Value = 3500
$ws.Range("M19").Value = -294.66666
$ws.Range("N19").Value = -3840
$ws.Range("H24").Value = 898.2857
$ws.Range("I24").Value = 464.66666
$ws.Range("J24").Value = 3500
$ws.Range("K24").Value = 464.66666
$ws.Range("L24").Value = 3500
$ws.Range("M24").Value = -294.66666
$ws.Range("N24").Value = -3840
$ws.Range("H29").Value = 5161.1333
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 5161.1333
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 5161.1333
$ws.Range("N29").Value = -5747.1333
$ws.Range("H35").Value = 655.55554
$ws.Range("I35").Value = 655.55554
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 655.55554
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -361.55554
$ws.Range("N35").Value = $null
$ws.Range("H122").Value = 64090.25
$ws.Range("I122").Value = 78511.08
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 235533.24
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -233083.24
$ws.Range("N122").Value = -9700

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 2797.4614
$ws.Range("I15").Value = 539.25
$ws.Range("J15").Value = 3801.111
$ws.Range("K15").Value = 1617.75
$ws.Range("L15").Value = 11403.333
$ws.Range("M15").Value = -1477.75
$ws.Range("N15").Value = -11683.333
$ws.Range("H122").Value = 10990
$ws.Range("I122").Value = 25474.666
$ws.Range("J122").Value = 2299.2
$ws.Range("K122").Value = 229271.994
$ws.Range("L122").Value = 20692.8
$ws.Range("M122").Value = -226821.994
$ws.Range("N122").Value = -25592.8
$ws.Range("H140").Value = 1130
$ws.Range("I140").Value = 1130
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 3390
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 1790
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15313.4
$ws.Range("I70").Value = 13423.25
$ws.Range("J70").Value = 22874
$ws.Range("K70").Value = 13423.25
$ws.Range("L70").Value = 22874
$ws.Range("M70").Value = -13153.25
$ws.Range("N70").Value = -23414
$ws.Range("H73").Value = 15313.4
$ws.Range("I73").Value = 13423.25
$ws.Range("J73").Value = 22874
$ws.Range("K73").Value = 13423.25
$ws.Range("L73").Value = 22874
$ws.Range("M73").Value = -12487.25
$ws.Range("N73").Value = -24746
$ws.Range("H97").Value = 855.82355
$ws.Range("I97").Value = 835.46155
$ws.Range("J97").Value = 922
$ws.Range("K97").Value = 835.46155
$ws.Range("L97").Value = 922
$ws.Range("M97").Value = -339.46155
$ws.Range("N97").Value = -1914
$ws.Range("H122").Value = 5455.294
$ws.Range("I122").Value = 5116
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 15348
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -12898
$ws.Range("N122").Value = -28900
$ws.Range("H141").Value = 193484.6
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 193484.6
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 193484.6
$ws.Range("N141").Value = -203844.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4194.3555
$ws.Range("I46").Value = 1703.5385
$ws.Range("J46").Value = 5206.25
$ws.Range("K46").Value = 1703.5385
$ws.Range("L46").Value = 5206.25
$ws.Range("M46").Value = -1515.5385
$ws.Range("N46").Value = -5582.25
$ws.Range("H61").Value = 2813.2307
$ws.Range("I61").Value = 3041.261
$ws.Range("J61").Value = 1065
$ws.Range("K61").Value = 3041.261
$ws.Range("L61").Value = 1065
$ws.Range("M61").Value = -2839.261
$ws.Range("N61").Value = -1469
$ws.Range("H100").Value = 62507500
$ws.Range("I100").Value = 125005000
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 125005000
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -125004459
$ws.Range("N100").Value = -11082
$ws.Range("H113").Value = 2813.2307
$ws.Range("I113").Value = 3041.261
$ws.Range("J113").Value = 1065
$ws.Range("K113").Value = 3041.261
$ws.Range("L113").Value = 1065
$ws.Range("M113").Value = -871.261
$ws.Range("N113").Value = -5405
$ws.Range("H122").Value = 5191.4443
$ws.Range("I122").Value = 2953.8333
$ws.Range("J122").Value = 9666.666999999999
$ws.Range("K122").Value = 8861.499899999999
$ws.Range("L122").Value = 29000.001
$ws.Range("M122").Value = -6411.499899999999
$ws.Range("N122").Value = -33900.001
$ws.Range("H132").Value = 7339.5835
$ws.Range("I132").Value = 6906
$ws.Range("J132").Value = 8206.75
$ws.Range("K132").Value = 20718
$ws.Range("L132").Value = 24620.25
$ws.Range("M132").Value = -18188
$ws.Range("N132").Value = -29680.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 500247.5
$ws.Range("I8").Value = 500247.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 500247.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -500107.5
$ws.Range("N8").Value = $null
$ws.Range("H62").Value = 30677.2
$ws.Range("I62").Value = 67497.5
$ws.Range("J62").Value = 6130.3335
$ws.Range("K62").Value = 67497.5
$ws.Range("L62").Value = 6130.3335
$ws.Range("M62").Value = -66873.5
$ws.Range("N62").Value = -7378.3335
$ws.Range("H65").Value = 30677.2
$ws.Range("I65").Value = 67497.5
$ws.Range("J65").Value = 6130.3335
$ws.Range("K65").Value = 337487.5
$ws.Range("L65").Value = 30651.6675
$ws.Range("M65").Value = -334367.5
$ws.Range("N65").Value = -36891.6675
$ws.Range("H132").Value = 5137.4683
$ws.Range("I132").Value = 3693.7144
$ws.Range("J132").Value = 7265.1055
$ws.Range("K132").Value = 11081.1432
$ws.Range("L132").Value = 21795.3165
$ws.Range("M132").Value = -8551.143199999999
$ws.Range("N132").Value = -26855.3165
